$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132261037826538
$ws.Range("B1").Value = 2.104254245758057
$ws.Range("C1").Value = 10.07627773284912
$ws.Range("D1").Value = 2.509932279586792
$ws.Range("E1").Value = 1.294487714767456
